$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps numeric-looking values such as
# "1.001" or "0.06590" as plain text, matching the source data feed's
# formatting instead of letting Excel coerce them into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.640.77"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "1.830.33"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "316.43"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").Value = "0.4001"
$ws.Range("E8").Value = "  +5.71%  "
$ws.Range("D9").Value = "0.07786"
$ws.Range("E9").Value = "  +4.02%  "
$ws.Range("D10").Value = "1.123"
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("D11").Value = "42.04"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "21.28"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").Value = "6.338"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "7.606"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "1.826.35"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "93.24"
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("D18").Value = "0.00001096"
$ws.Range("E18").Value = "  +2.78%  "
$ws.Range("D19").Value = "0.06590"
$ws.Range("D20").Value = "17.85"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "6.109"
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("D23").Value = "28.632.82"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "2.232"
$ws.Range("E25").Value = "  +6.74%  "
$ws.Range("D26").Value = "20.88"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").Value = "156.67"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").Value = "2.036.64"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Value = "2.435"
$ws.Range("E29").Value = "  +4.00%  "
$ws.Range("D30").Value = "125.77"
$ws.Range("E30").Value = "  +2.88%  "
$ws.Range("D31").Value = "1.167"
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("D32").Value = "0.1129"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").Value = "5.773"
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("D34").Value = "3.662"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "0.07381"
$ws.Range("E35").Value = "  +5.53%  "
$ws.Range("D36").Value = "0.2280"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("D37").Value = "0.02360"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("D38").Value = "8.953"
$ws.Range("E38").Value = "  +5.66%  "
$ws.Range("D39").Value = "5.223"
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("D40").Value = "11.43"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("D41").Value = "0.6317"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").Value = "1.200"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").Value = "13.61"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").Value = "0.5961"
$ws.Range("E46").Value = "  +3.10%  "
$ws.Range("D47").Value = "3.715"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").Value = "125.91"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").Value = "2.006"
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("D50").Value = "1.196"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "0.06972"
$ws.Range("E51").Value = "  +2.15%  "

